# Insert a new price record row before row 384 (Feria Lagunitas de Puerto
# Montt, Pina, Segunda / Ecuador), shifting the existing rows 384-458 down
# to 385-459, then populate the newly inserted row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(384).EntireRow.Insert()

$ws.Cells.Item(384, 1).Value  = 4
$ws.Cells.Item(384, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(384, 3).Value  = "Los Lagos"
$ws.Cells.Item(384, 4).Value  = 45173
$ws.Cells.Item(384, 5).Value  = 10
$ws.Cells.Item(384, 6).Value  = "Fruta"
$ws.Cells.Item(384, 7).Value  = 100108
$ws.Cells.Item(384, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(384, 9).Value  = 100108005
$ws.Cells.Item(384, 10).Value = "Pi$([char]0x00F1)a"
$ws.Cells.Item(384, 11).Value = "Caramelo"
$ws.Cells.Item(384, 12).Value = "Segunda"
$ws.Cells.Item(384, 13).Value = 60
$ws.Cells.Item(384, 14).Value = 25000
$ws.Cells.Item(384, 15).Value = 25000
$ws.Cells.Item(384, 16).Value = 25000
$ws.Cells.Item(384, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(384, 18).Value = "Ecuador"
$ws.Cells.Item(384, 19).Value = 1786
$ws.Cells.Item(384, 20).Value = 14
